$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25 / Row 26: shift existing B/C into A/C, insert new B text ---
$ws.Range("A25").Value2 = "Cat-Rules-1"
$ws.Range("B25").Value2 = "Cannot modify menu if it is in use"

$ws.Range("A26").Value2 = "Security-5"
$ws.Range("B26").Value2 = "Can relog in"

# --- Row 27 (new) ---
$ws.Range("A27").Value2 = "Security-4"
$ws.Range("B27").Value2 = "The password requirements are tested"
$ws.Range("C27").Value2 = "Each rule broken rule should return a unique binary value"

# --- Row 28 / Row 29 (new) ---
$ws.Range("A28").Value2 = "User-Factory-2"
$ws.Range("A29").Value2 = "User-Service-6"

$ws.Range("B28").Value2 = "Get Result<Token>"
$ws.Range("B29").Value2 = "User location has been changed."

$ws.Range("C29").Value2 = "Null and white-space strings should not change the location"

$ws.Range("C28").Value2 = "Should return either SuccessResullt or BadRequest"
$ws.Range("D28").Value2 = "SuccessResult if data is valid else BadRequest"
$ws.Range("E28").Value2 = "All paths tested"

# --- Row 30 (new) ---
$ws.Range("A30").Value2 = "Cat-Service-16"
$ws.Range("B30").Value2 = "Customer location has been changed."
$ws.Range("C30").Value2 = "Null and white-space strings should not change the location"

# --- Row 31 (new) ---
$ws.Range("A31").Value2 = "User-Endpoint-1"
$ws.Range("C31").Value2 = "200 with Auth data if valid user else 401"
$ws.Range("B31").Value2 = "Attemp to log in"

# --- Row 23 / Row 24: update Notes column text ---
$ws.Range("D23").Value2 = "Mock the communication interface. Requires being loggedin"
$ws.Range("D24").Value2 = "Mock the communication interface. Requires being loggedin"

# --- Row 32 (new) ---
$ws.Range("A32").Value2 = "Result-1-a"
$ws.Range("B32").Value2 = "Mapping Result<T>/Result to correct HTTPResponse"
$ws.Range("C32").Value2 = "Correct mapping. 200 should have data, all non-200 and non-204 should carry the errors"

# --- Row 33 (new) ---
$ws.Range("A33").Value2 = "User-Factory-1"
$ws.Range("B33").Value2 = "Get Result<User>"
$ws.Range("C33").Value2 = "Should return either SuccessResullt or BadRequest"

# Wrap text for A29
$ws.Range("A29").WrapText = $true

# Fix selection / view
$ws.Range("B31").Select()

$wb.Save()
